$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 351266
$ws.Range("E2").Value = 31412
$ws.Range("F2").Value = 30706
$ws.Range("G2").Value = 46590
$ws.Range("H2").Value = 33925
$ws.Range("I2").Value = 34224
$ws.Range("J2").Value = -299
$ws.Range("K2").Value = 391119
$ws.Range("L2").Value = 158258
$ws.Range("M2").Value = 232861
$ws.Range("N2").Value = 231627
$ws.Range("O2").Value = 1234
$ws.Range("P2").Value = 4911
$ws.Range("Q2").Value = 8466
$ws.Range("R2").Value = -1866
$ws.Range("S2").Value = -2223
$ws.Range("T2").Value = 10695
$ws.Range("U2").Value = -2229
$ws.Range("V2").Value = 28604
$ws.Range("W2").Value = 8.94
$ws.Range("X2").Value = 9.66
$ws.Range("Y2").Value = 15.82
$ws.Range("Z2").Value = 9.23
$ws.Range("AA2").Value = 67.95999999999999
$ws.Range("AB2").Value = 4767.1
$ws.Range("AC2").Value = 35157
$ws.Range("AD2").Value = 6.71
$ws.Range("AE2").Value = 242491
$ws.Range("AF2").Value = 0.97
$ws.Range("AG2").Value = 3000
$ws.Range("AH2").Value = 1.27
$ws.Range("AI2").Value = 8.369999999999999
$ws.Range("AJ2").Value = 97343863
$ws.Range("D3").Value = 360197
$ws.Range("E3").Value = 29346
$ws.Range("F3").Value = 29346
$ws.Range("G3").Value = 42127
$ws.Range("H3").Value = 30400
$ws.Range("I3").Value = 30554
$ws.Range("J3").Value = -154
$ws.Range("K3").Value = 377748
$ws.Range("L3").Value = 120986
$ws.Range("M3").Value = 256762
$ws.Range("N3").Value = 256218
$ws.Range("O3").Value = 545
$ws.Range("P3").Value = 4911
$ws.Range("Q3").Value = 43519
$ws.Range("R3").Value = -49486
$ws.Range("S3").Value = 1698
$ws.Range("T3").Value = 40798
$ws.Range("U3").Value = 2720
$ws.Range("V3").Value = 32218
$ws.Range("W3").Value = 8.15
$ws.Range("X3").Value = 8.44
$ws.Range("Y3").Value = 12.53
$ws.Range("Z3").Value = 7.91
$ws.Range("AA3").Value = 47.12
$ws.Range("AB3").Value = 5335.93
$ws.Range("AC3").Value = 31387
$ws.Range("AD3").Value = 7.85
$ws.Range("AE3").Value = 270935
$ws.Range("AF3").Value = 0.91
$ws.Range("AG3").Value = 3500
$ws.Range("AH3").Value = 1.42
$ws.Range("AI3").Value = 10.83
$ws.Range("AJ3").Value = 97343863
$ws.Range("D4").Value = 382617
$ws.Range("E4").Value = 29047
$ws.Range("F4").Value = 29047
$ws.Range("G4").Value = 41112
$ws.Range("H4").Value = 30473
$ws.Range("I4").Value = 30378
$ws.Range("J4").Value = 95
$ws.Range("K4").Value = 417116
$ws.Range("L4").Value = 131536
$ws.Range("M4").Value = 285580
$ws.Range("N4").Value = 284945
$ws.Range("O4").Value = 635
$ws.Range("P4").Value = 4911
$ws.Range("Q4").Value = 20178
$ws.Range("R4").Value = -21245
$ws.Range("S4").Value = -3363
$ws.Range("T4").Value = 12961
$ws.Range("U4").Value = 7217
$ws.Range("V4").Value = 32907
$ws.Range("W4").Value = 7.59
$ws.Range("X4").Value = 7.96
$ws.Range("Y4").Value = 11.23
$ws.Range("Z4").Value = 7.67
$ws.Range("AA4").Value = 46.06
$ws.Range("AB4").Value = 5890.06
$ws.Range("AC4").Value = 31205
$ws.Range("AD4").Value = 8.460000000000001
$ws.Range("AE4").Value = 300878
$ws.Range("AF4").Value = 0.88
$ws.Range("AG4").Value = 3500
$ws.Range("AH4").Value = 1.33
$ws.Range("AI4").Value = 10.91
$ws.Range("AJ4").Value = 97343863
$ws.Range("D5").Value = 351446
$ws.Range("E5").Value = 20249
$ws.Range("F5").Value = 20249
$ws.Range("G5").Value = 27344
$ws.Range("H5").Value = 15577
$ws.Range("I5").Value = 15682
$ws.Range("J5").Value = -104
$ws.Range("K5").Value = 417368
$ws.Range("L5").Value = 123779
$ws.Range("M5").Value = 293590
$ws.Range("N5").Value = 292954
$ws.Range("O5").Value = 636
$ws.Range("P5").Value = 4911
$ws.Range("Q5").Value = 19487
$ws.Range("R5").Value = -10655
$ws.Range("S5").Value = -3955
$ws.Range("T5").Value = 6769
$ws.Range("U5").Value = 12719
$ws.Range("V5").Value = 30667
$ws.Range("W5").Value = 5.76
$ws.Range("X5").Value = 4.43
$ws.Range("Y5").Value = 5.43
$ws.Range("Z5").Value = 3.73
$ws.Range("AA5").Value = 42.16
$ws.Range("AB5").Value = 6146.92
$ws.Range("AC5").Value = 16109
$ws.Range("AD5").Value = 16.33
$ws.Range("AE5").Value = 309334
$ws.Range("AF5").Value = 0.85
$ws.Range("AG5").Value = 3500
$ws.Range("AH5").Value = 1.33
$ws.Range("AI5").Value = 21.14
$ws.Range("AJ5").Value = 97343863
$ws.Range("D6").Value = 351492
$ws.Range("E6").Value = 20250
$ws.Range("F6").Value = 20250
$ws.Range("G6").Value = 24749
$ws.Range("H6").Value = 18882
$ws.Range("I6").Value = 18888
$ws.Range("K6").Value = 430711
$ws.Range("L6").Value = 123677
$ws.Range("M6").Value = 307034
$ws.Range("N6").Value = 306305
$ws.Range("P6").Value = 4911
$ws.Range("Q6").Value = 16101
$ws.Range("R6").Value = -9440
$ws.Range("S6").Value = -7206
$ws.Range("T6").Value = 5314
$ws.Range("U6").Value = 10787
$ws.Range("V6").Value = 27011
$ws.Range("W6").Value = 5.76
$ws.Range("X6").Value = 5.37
$ws.Range("Y6").Value = 6.3
$ws.Range("Z6").Value = 4.45
$ws.Range("AA6").Value = 40.28
$ws.Range("AB6").Value = 6412.49
$ws.Range("AC6").Value = 19403
$ws.Range("AD6").Value = 9.789999999999999
$ws.Range("AE6").Value = 323432
$ws.Range("AF6").Value = 0.59
$ws.Range("AG6").Value = 4000
$ws.Range("AH6").Value = 2.11
$ws.Range("AI6").Value = 20.06
$ws.Range("AJ6").Value = 97343863
$ws.Range("D7").Value = 376390
$ws.Range("E7").Value = 23596
$ws.Range("G7").Value = 32083
$ws.Range("H7").Value = 23619
$ws.Range("I7").Value = 23572
$ws.Range("K7").Value = 459635
$ws.Range("L7").Value = 130786
$ws.Range("M7").Value = 328837
$ws.Range("N7").Value = 327627
$ws.Range("P7").Value = 4905
$ws.Range("Q7").Value = 27837
$ws.Range("R7").Value = -12844
$ws.Range("S7").Value = -4322
$ws.Range("T7").Value = 7195
$ws.Range("U7").Value = 18409
$ws.Range("W7").Value = 6.27
$ws.Range("X7").Value = 6.28
$ws.Range("Y7").Value = 7.44
$ws.Range("Z7").Value = 5.31
$ws.Range("AA7").Value = 39.77
$ws.Range("AC7").Value = 24561
$ws.Range("AD7").Value = 10.24
$ws.Range("AE7").Value = 350623
$ws.Range("AF7").Value = 0.72
$ws.Range("AG7").Value = 4413
$ws.Range("AH7").Value = 1.75
$ws.Range("AI7").Value = 17.84
$ws.Range("D8").Value = 404238
$ws.Range("E8").Value = 26602
$ws.Range("G8").Value = 37123
$ws.Range("H8").Value = 27533
$ws.Range("I8").Value = 27485
$ws.Range("K8").Value = 487880
$ws.Range("L8").Value = 137969
$ws.Range("M8").Value = 349915
$ws.Range("N8").Value = 349383
$ws.Range("P8").Value = 4904
$ws.Range("Q8").Value = 26152
$ws.Range("R8").Value = -15452
$ws.Range("S8").Value = -2602
$ws.Range("T8").Value = 9491
$ws.Range("U8").Value = 17254
$ws.Range("W8").Value = 6.58
$ws.Range("X8").Value = 6.81
$ws.Range("Y8").Value = 8.119999999999999
$ws.Range("Z8").Value = 5.81
$ws.Range("AA8").Value = 39.43
$ws.Range("AC8").Value = 28907
$ws.Range("AD8").Value = 7.94
$ws.Range("AE8").Value = 374918
$ws.Range("AF8").Value = 0.61
$ws.Range("AG8").Value = 4802
$ws.Range("AH8").Value = 2.09
$ws.Range("AI8").Value = 16.65
$ws.Range("D9").Value = 432223
$ws.Range("E9").Value = 29396
$ws.Range("G9").Value = 40915
$ws.Range("H9").Value = 30347
$ws.Range("I9").Value = 30289
$ws.Range("K9").Value = 520468
$ws.Range("L9").Value = 145322
$ws.Range("M9").Value = 375152
$ws.Range("N9").Value = 374578
$ws.Range("P9").Value = 4904
$ws.Range("Q9").Value = 27102
$ws.Range("R9").Value = -15206
$ws.Range("S9").Value = -2037
$ws.Range("T9").Value = 10033
$ws.Range("U9").Value = 17333
$ws.Range("W9").Value = 6.8
$ws.Range("X9").Value = 7.02
$ws.Range("Y9").Value = 8.369999999999999
$ws.Range("Z9").Value = 6.02
$ws.Range("AA9").Value = 38.74
$ws.Range("AC9").Value = 31863
$ws.Range("AD9").Value = 7.2
$ws.Range("AE9").Value = 401955
$ws.Range("AF9").Value = 0.57
$ws.Range("AG9").Value = 5026
$ws.Range("AH9").Value = 2.19
$ws.Range("AI9").Value = 15.82
